$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new rows ("line7", "line8") are inserted right after the existing
# "line6" row (row 7) and before the old "extr1" row (old row 8), pushing
# the old rows 8-15 ("extr1".."extr8") down to rows 10-17. Shift the
# block down two rows by copying it (this also extends the used range /
# dimension to row 17 and keeps each cell's existing formatting).
$ws.Range("A8:E15").Copy($ws.Range("A10:E15"))

# Full target data for rows 8-17 (A=index, B=name, C=from_bus, D=to_bus, E=in_service).
$data = @(
    @(8,  6,  "line7", 14, 11, $true),
    @(9,  7,  "line8", 16, 9,  $true),
    @(10, 8,  "extr1", 5,  12, $true),
    @(11, 9,  "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $true),
    @(13, 11, "extr4", 7,  8,  $true),
    @(14, 12, "extr5", 9,  11, $true),
    @(15, 13, "extr6", 7,  11, $true),
    @(16, 14, "extr7", 5,  7,  $true),
    @(17, 15, "extr8", 8,  5,  $true)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
}
